# "Add GradeSheet, simplify Indicator, Criterion and Evaluation"
#
# On the "gc_rhum" grading sheet, rows 16-18 (the three indicators under
# criterion C2) were re-graded using the raw-points column (C) instead of
# the manual percentage column (B): column B is cleared and column C now
# holds the indicator's point total (13/14, 13/13, 13/13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gc_rhum")

# Indicator "Cordage enroulé en huit impeccable (14 pts)"
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 13

# Indicator "Pavillon solidement fixé au mât (13 pts)"
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 13

# Indicator "Drisse tendue sans nœud lâche (13 pts)"
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = 13

# Reflect the user's last click on the sheet.
$ws.Activate()
$ws.Range("C16").Select()
